# Regenerate the "K" (strikeouts) column (G) in the save_data sheet for
# kennedy_ian 2022 using the recomputed K values (previously a placeholder
# "Strike#" figure). Data rows run from row 2 (game index 0) to row 63
# (game index 61); column G is the 7th column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for G2:G63, in row order.
$kValues = @(
    1,1,1,1,0,
    0,0,0,0,1,
    2,0,1,0,1,
    0,1,1,1,3,
    0,2,0,0,2,
    0,1,1,1,0,
    2,1,0,2,1,
    1,0,0,1,2,
    2,1,1,0,0,
    0,0,1,1,1,
    1,2,1,0,1,
    0,0,0,0,0,
    0,0
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
